# Template_ACI _STE_ADDAX.xlsx — "Summary Features" touch-up
#
# 1) Clean up the "DUE_AMT " column header (trailing space removed).
# 2) Re-apply cell protection on the header cell D1 (Format Cells > Protection),
#    which is what produces the new/updated cell style entry in the workbook.
# 3) Leave the selection where the author last left it (G6) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the shared-string header text for column D.
$ws.Range("D1").Value = "DUE_AMT"

# 2) Touch the Protection settings of the header cell so the workbook records
#    an explicit "locked" cell format for it (adds applyProtection to its style).
$ws.Range("D1").Locked = $true

# 3) Restore the last active selection used when the file was saved.
$ws.Range("G6").Select() | Out-Null
